# Prefix every "StepN ..." / short command name in column A (rows 2..N)
# of each protocol worksheet with that worksheet's own name, e.g.
#   "Step4 Seed" on sheet "free1"  ->  "free1 Step4 Seed"
#
# This mirrors the commit: "fix: unique command names in XLSX - prefix
# protocol name to each step"

$wb = $excel.ActiveWorkbook

# Sheets whose column-A command names need the sheet name prefixed on.
$sheetNames = @(
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol",
    "dickpic",
    "boosters",
    "price1", "price2",
    "discount1", "discount2"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $used = $ws.UsedRange
    $firstRow = $used.Row
    $lastRow = $firstRow + $used.Rows.Count - 1

    # Row 1 is the header ("Name", "Text", "Note", "*Guidelines"); data
    # rows start at row 2 and run through the last used row.
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2

        if ($null -ne $current -and $current -ne "") {
            $prefix = $sheetName + " "
            if (-not $current.StartsWith($prefix)) {
                $cell.Value2 = $prefix + $current
            }
        }
    }
}
